$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 90.353905
$ws.Range("N2").Value = 180.70781
$ws.Range("O2").Value = 0.3131638580342592
$ws.Range("P2").Value = 0.2965570309229201
$ws.Range("Q2").Value = 15.65239849673833
$ws.Range("R2").Value = 93.91439098043
$ws.Range("S2").Value = 0.3131638580342592
$ws.Range("T2").Value = 0.2965570309229201

# Row 3
$ws.Range("O3").Value = 0.001245676287098259
$ws.Range("P3").Value = 0.001769428615638018
$ws.Range("S3").Value = 0.001245676287098259
$ws.Range("T3").Value = 0.001769428615638018

# Row 4
$ws.Range("M4").Value = 10.80810533333333
$ws.Range("N4").Value = 32.424316
$ws.Range("O4").Value = 0.03746056093787335
$ws.Range("P4").Value = 0.05321108635352579
$ws.Range("Q4").Value = 1.872334922016444
$ws.Range("R4").Value = 16.851014298148
$ws.Range("S4").Value = 0.03746056093787335
$ws.Range("T4").Value = 0.05321108635352579

# Row 5
$ws.Range("M5").Value = 165.852196
$ws.Range("N5").Value = 331.704392
$ws.Range("O5").Value = 0.574838614477306
$ws.Range("P5").Value = 0.5443553858331436
$ws.Range("Q5").Value = 28.73129460592934
$ws.Range("R5").Value = 172.387767635576
$ws.Range("S5").Value = 0.574838614477306
$ws.Range("T5").Value = 0.5443553858331436

# Row 6
$ws.Range("M6").Value = 5.448456
$ws.Range("N6").Value = 16.345368
$ws.Range("O6").Value = 0.01888418105769649
$ws.Range("P6").Value = 0.02682415222353981
$ws.Range("Q6").Value = 0.9438596428560001
$ws.Range("R6").Value = 8.494736785704001
$ws.Range("S6").Value = 0.01888418105769649
$ws.Range("T6").Value = 0.02682415222353981

# Row 7
$ws.Range("M7").Value = 15.69751633333333
$ws.Range("N7").Value = 47.09254900000001
$ws.Range("O7").Value = 0.05440710920576665
$ws.Range("P7").Value = 0.07728291605123282
$ws.Range("Q7").Value = 2.719348776994112
$ws.Range("R7").Value = 24.474138992947
$ws.Range("S7").Value = 0.05440710920576665
$ws.Range("T7").Value = 0.07728291605123282
